# ------------------------------------------------------------------
# Updates in project framework
#  - insert a new blank "Sheet2" tab after "contact"
#  - insert new "CampaignModule" and "OpportunitiesModule" tabs after
#    "ProductModule", populated with test-case rows
#  - tweak a couple of view selections on existing tabs
#  - make "OpportunitiesModule" the active tab
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. new blank "Sheet2" tab, placed right after "contact" ------
$contact = $wb.Worksheets.Item("contact")
$sheet2 = $wb.Worksheets.Add($null, $contact)
$sheet2.Name = "Sheet2"

# ---- 2. new "CampaignModule" tab, placed right after "ProductModule"
$productModule = $wb.Worksheets.Item("ProductModule")
$campaign = $wb.Worksheets.Add($null, $productModule)
$campaign.Name = "CampaignModule"

$campaign.Range("A1").Value = "TC_ID"
$campaign.Range("B1").Value = "TC_NAME"
$campaign.Range("C1").Value = "CampaignName"
$campaign.Range("D1").Value = "EventName"

$campaign.Range("A2").Value = "TC_01"
$campaign.Range("B2").Value = "Create campaign with events"
$campaign.Range("C2").Value = "Greenliving"
$campaign.Range("D2").Value = "Breezy Brights"

$campaign.Range("A3").Value = "TC_02"
$campaign.Range("B3").Value = "Search campaign and delete"
$campaign.Range("C3").Value = "Greenliving"
$campaign.Range("D3").Value = "Breezy Brights"

# formatting: header row gets the same light header-fill used elsewhere
# in the workbook, body uses text ("@") number format
$campaign.Range("A1:D3").NumberFormat = "@"
$campaign.Range("A1:D1").Interior.Color = 16250605
$campaign.Columns.Item(1).ColumnWidth = 9.140625
$campaign.Columns.Item(2).ColumnWidth = 27.140625
$campaign.Columns.Item(3).ColumnWidth = 15.140625
$campaign.Columns.Item(4).ColumnWidth = 13.7109375
$campaign.Range("D20").Select()

# ---- 3. new "OpportunitiesModule" tab, placed right after "CampaignModule"
$opportunities = $wb.Worksheets.Add($null, $campaign)
$opportunities.Name = "OpportunitiesModule"

$opportunities.Range("A1").Value = "TC_ID"
$opportunities.Range("B1").Value = "TC_NAME"
$opportunities.Range("C1").Value = "OrganizationName"
$opportunities.Range("D1").Value = "OpportunityName"
$opportunities.Range("E1").Value = "Expected_Close_Date"

$opportunities.Range("A2").Value = "TC_01"
$opportunities.Range("B2").Value = "Create opportunity with organization name and assign to group"
$opportunities.Range("C2").Value = "PIMCHA"
$opportunities.Range("D2").Value = "HempCann"
$opportunities.Range("E2").Value = 45527

$opportunities.Range("A3").Value = "TC_02"
$opportunities.Range("B3").Value = "Create opportunity with organization name with back date from current date"
$opportunities.Range("C3").Value = "PIMCHA"
$opportunities.Range("D3").Value = "HempCann"
$opportunities.Range("E3").Value = 44795

$opportunities.Range("A1:E3").NumberFormat = "@"
$opportunities.Range("A1:G1").Interior.Color = 12566463
$opportunities.Columns.Item(2).ColumnWidth = 70.42578125
$opportunities.Columns.Item(3).ColumnWidth = 17.85546875
$opportunities.Columns.Item(4).ColumnWidth = 17.42578125
$opportunities.Columns.Item(5).ColumnWidth = 20.42578125
$opportunities.Range("B15").Select()

# ---- 4. minor view-selection tweaks on existing tabs --------------
$org = $wb.Worksheets.Item("org")
$org.Range("B14").Select()

$productModule.Range("A1:D1").Select()

# ---- 5. make "OpportunitiesModule" the active/visible tab ---------
$opportunities.Activate()
